$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# "Last updated" banner timestamp bump (07:06 -> 08:23)
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 08:23"

# Refreshed per-country COVID figures. Three country rows (190-192 and
# 196-197) also had their labels re-sorted alphabetically, which shifted
# which row each country now sits on.

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 795605
$ws.Range("C6").Value = 763
$ws.Range("D6").Value = 496048
$ws.Range("E6").Value = 277925
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = 21632

# Row 71: Uzbekistan
$ws.Range("A71").Value = "Uzbekistan"
$ws.Range("B71").Value = 11723
$ws.Range("C71").Value = 159
$ws.Range("D71").Value = 7287
$ws.Range("E71").Value = 4384
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 52

# Row 77: El Salvador
$ws.Range("A77").Value = "El Salvador"
$ws.Range("B77").Value = 8844
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 5341
$ws.Range("E77").Value = 3254
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 6
$ws.Range("H77").Value = 249

# Row 102: Tailandia
$ws.Range("A102").Value = "Tailandia"
$ws.Range("B102").Value = 3202
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 3087
$ws.Range("E102").Value = 57
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 58

# Row 112: Sri Lanka
$ws.Range("A112").Value = "Sri Lanka"
$ws.Range("B112").Value = 2350
$ws.Range("C112").Value = 196
$ws.Range("D112").Value = 1979
$ws.Range("E112").Value = 360
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 11

# Row 116: Malaui
$ws.Range("A116").Value = "Malaui"
$ws.Range("B116").Value = 1984
$ws.Range("C116").Value = 42
$ws.Range("D116").Value = 369
$ws.Range("E116").Value = 1590
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 25

# Row 143: Georgia
$ws.Range("A143").Value = "Georgia"
$ws.Range("B143").Value = 973
$ws.Range("C143").Value = 5
$ws.Range("D143").Value = 846
$ws.Range("E143").Value = 112
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 15

# Row 145: Zimbabue
$ws.Range("A145").Value = "Zimbabue"
$ws.Range("B145").Value = 926
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 306
$ws.Range("E145").Value = 608
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 12

# Row 190: Islas Turcas y Caicos (was Gambia on this row)
$ws.Range("A190").Value = "Islas Turcas y Caicos"
$ws.Range("B190").Value = 66
$ws.Range("C190").Value = 11
$ws.Range("D190").Value = 11
$ws.Range("E190").Value = 53
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 2

# Row 191: Gambia (was Polinesia Francesa on this row)
$ws.Range("A191").Value = "Gambia"
$ws.Range("B191").Value = 64
$ws.Range("C191").Value = 1
$ws.Range("D191").Value = 34
$ws.Range("E191").Value = 27
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 3

# Row 192: Polinesia Francesa (was Islas Turcas y Caicos on this row)
$ws.Range("A192").Value = "Polinesia Francesa"
$ws.Range("B192").Value = 62
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 60
$ws.Range("E192").Value = 2
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 0

# Row 196: Belice (was Guam on this row)
$ws.Range("A196").Value = "Belice"
$ws.Range("B196").Value = 33
$ws.Range("C196").Value = 3
$ws.Range("D196").Value = 20
$ws.Range("E196").Value = 11
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 2

# Row 197: Guam (was Belice on this row)
$ws.Range("A197").Value = "Guam"
$ws.Range("B197").Value = 32
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 0
$ws.Range("E197").Value = 31
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 1
